$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare styles for the new cells first (no text values yet, so shared-string
#     order below is driven purely by the Value assignments that follow) ---

# Row 8 (SlNo 6, continues the 30-Oct block)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Row 9: blank separator row, only B9 carries the date number-format style
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Rows 10-13 (31-Oct block)
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)

# --- Now assign the actual values/text, in the same order the strings were
#     first introduced so the rebuilt sharedStrings table lines up ---

# Row 3: status Partial -> In Progress
$ws.Range("D3").Value = "In Progress"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 45229
$ws.Range("C8").Value = "Mock Test Planning"
$ws.Range("D8").Value = "In Progress"
$ws.Range("E8").Value = "Called Waqar Younis and Manpreet of College Doors but they did n't picked up the call"

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 45230
$ws.Range("C10").Value = "Connected with Waqar Younis of College Doors about Mock Exam"
$ws.Range("D10").Value = "Done"

# Row 12 text (entered before row 11 text, matching the source order)
$ws.Range("C12").Value = "Talked with Mr Subroto Ghosh for practice test"

# Row 11
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 45230
$ws.Range("C11").Value = "Sayan took my AI class for 2 hours"
$ws.Range("D11").Value = "Done"

# Row 3 remarks (updated text)
$ws.Range("E3").Value = "Called Mr Subroto Ghosh and he said he will update the note details in night"

# Row 12 (remaining cells)
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 45230
$ws.Range("D12").Value = "Done"
$ws.Range("E12").Value = "For JILETE we need to purchase question bank and for Prottoy we need to talk with his father which test he wants to give board level or joint"

# Row 13
$ws.Range("B13").Value = 45230
$ws.Range("C13").Value = "Talked with Mr Amitabha Kairali for Prottoy's monthly test"
$ws.Range("D13").Value = "Done"
$ws.Range("E13").Value = "As Prottoy will give JEE and JEE WB exams, hence he will give practice test on this exams."

# Row 7: add missing Status "Done"
$ws.Range("D7").Value = "Done"

# Row 9: clear the value so only the style (date number format) remains
$ws.Range("B9").ClearContents()

# --- Column widths (col D new, col E widened) ---
$ws.Columns.Item(4).ColumnWidth = 9
$ws.Columns.Item(5).ColumnWidth = 100.4

# --- Sheet view: selection moves to E14 after data entry ---
[void]$ws.Range("E14").Select()
